# Daily attendance processing - 2026-01-21 23:38:20
# Swap the order of "Recorded By" entries that contain both the user email
# and "System" from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $v = $cell.Value()
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
